$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B slightly (from 14.42578125 to 15.42578125 character-units)
$ws.Range("A:B").ColumnWidth = 14.6

# Update the values in A1:B5 with new data
$ws.Range("A1").Value = -0.0077319541145460515
$ws.Range("B1").Value = -0.0073151622617322852

$ws.Range("A2").Value = -0.035113097906611136
$ws.Range("B2").Value = -0.045131186603830108

$ws.Range("A3").Value = -0.013595168236939547
$ws.Range("B3").Value = -0.025377727842637879

$ws.Range("A4").Value = -0.023474155373074377
$ws.Range("B4").Value = -0.022879431690503188

$ws.Range("A5").Value = -0.060146151198309898
$ws.Range("B5").Value = -0.060125882916858303
